# datasetOverview: add two new summary columns (mnRng / seRng) after the
# existing data, and switch the sheet view to a frozen-header/frozen-ID
# layout (freeze first 4 columns + header row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells -------------------------------------------------
# Write AX1 ("seRng") before AW1 ("mnRng") so the workbook's shared-string
# table gets the two new strings appended in the same order as the source
# edit (seRng, then mnRng).
$ws.Range("AX1").Value = "seRng"
$ws.Range("AW1").Value = "mnRng"

# --- New data columns ---------------------------------------------------
# Rows 7, 12, 19 and 20 have no range-size data for this pair of columns
# (same rows that are missing values further up the row, ending in #N/A),
# so they are intentionally skipped - only rows with data get AW/AX values.
$ws.Range("AW2").Value = 924.63
$ws.Range("AX2").Value = 57.04

$ws.Range("AW3").Value = 951.24
$ws.Range("AX3").Value = 64

$ws.Range("AW4").Value = 767.38
$ws.Range("AX4").Value = 52.04

$ws.Range("AW5").Value = 365.43
$ws.Range("AX5").Value = 34.83

$ws.Range("AW6").Value = 972.4
$ws.Range("AX6").Value = 50.75

$ws.Range("AW8").Value = 247.12
$ws.Range("AX8").Value = 47.61

$ws.Range("AW9").Value = 353.42
$ws.Range("AX9").Value = 26.39

$ws.Range("AW10").Value = 302.56
$ws.Range("AX10").Value = 34.16

$ws.Range("AW11").Value = 601.9
$ws.Range("AX11").Value = 64.21

$ws.Range("AW13").Value = 636.92
$ws.Range("AX13").Value = 67.1

$ws.Range("AW14").Value = 499.97
$ws.Range("AX14").Value = 69.64

$ws.Range("AW15").Value = 378.57
$ws.Range("AX15").Value = 77.47

$ws.Range("AW16").Value = 281.66
$ws.Range("AX16").Value = 20.45

$ws.Range("AW17").Value = 519.54
$ws.Range("AX17").Value = 52.59

$ws.Range("AW18").Value = 423.17
$ws.Range("AX18").Value = 22.27

$ws.Range("AW21").Value = 580.12
$ws.Range("AX21").Value = 30.31

# --- Sheet view: freeze panes --------------------------------------------
# Freeze the first 4 columns (A:D) and the header row (row 1). Selecting
# E2 (the cell just below/right of the freeze boundary) before turning on
# FreezePanes is what drives Excel to split at column D / row 1.
$ws.Range("E2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Leave the final selection on the last cell of the new data, matching
# the bottom-right pane's active cell/selection.
$ws.Range("AX21").Select()
